# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Recomputes the K value for each existing data row and writes it back in place,
# leaving every other column (date, TB, PC, dS0, dSF, IP, I0, IF) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> recalculated K value (column G), derived/regen'd from the
# underlying box-score source for this save (K replaces the old Strike# count).
$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 2
    9  = 2
    10 = 0
    11 = 1
    12 = 2
    13 = 1
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 0
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 2
    31 = 1
    32 = 1
    33 = 0
    34 = 1
    35 = 2
    36 = 2
    39 = 0
    40 = 1
    41 = 1
    42 = 2
    43 = 2
    44 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
